{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the 5 textual changes described by the diff:\n//   1. Title: \"Github Foundations Certification Training\"\n//      -> \"Hands-on AI-Assisted Programming Made Simple with GitHub Copilot\"\n//   2. Funding paragraph: course name + framework code update\n//   3. \"Performance gaps\" paragraph rewrite\n//   4. \"Why this course\" paragraph rewrite\n//   5. Date line: \"25 February 2025\" -> \"03 March 2025\"\n\nconst replacements = [\n  // NOTE: the longer sentence containing \"Github Foundations Certification\n  // Training\" as a substring must be replaced BEFORE the standalone title\n  // occurrence, otherwise a search for the short title also matches (and\n  // clobbers) the occurrence embedded inside this sentence.\n  {\n    find:\n      \"We are applying for WSQ funding support for this new course Github Foundations Certification Training according to Software Configuration ICT-DIT-3014-1.1 under Infocomm Technology Framework.\",\n    replace:\n      \"We are applying for WSQ funding support for this new course Hands-on AI-Assisted Programming Made Simple with GitHub Copilot according to Digital Technology Adoption and Innovation ACC-ICT-3004-1.1 under Infocomm Technology Framework.\",\n  },\n  {\n    find: \"Github Foundations Certification Training\",\n    replace: \"Hands-on AI-Assisted Programming Made Simple with GitHub Copilot\",\n  },\n  {\n    find:\n      \"Many software development teams face challenges in selecting and utilizing the right tools for integration and deployment. A lack of expertise in scripting and automation leads to manual processes that are time-consuming and prone to errors. Moreover, teams often struggle to properly diagnose configuration issues, resulting in prolonged debugging cycles and delayed project timelines.\",\n    replace:\n      \"One significant challenge is the slow adoption of new technologies and methodologies, hindering the ability to remain competitive. Teams may lack the expertise to effectively integrate emerging tools, which can lead to missed opportunities for automation and improved productivity. Legacy systems and a reluctance to change further compound these issues.\",\n  },\n  {\n    find:\n      \"The course covers a broad range of Git scripts and tools necessary for effectively integrating and deploying software. Participants will learn how to interpret configuration test results and identify the root causes of issues. This skillset enables them to proactively address problems and implement necessary modifications, which in turn leads to more reliable and efficient software releases.\",\n    replace:\n      \"This course directly addresses this by providing hands-on experience with cutting-edge AI programming tools. Participants will explore how these tools can streamline organizational coding processes, including using code completion and suggesting code snippets. The course is designed to help individuals stay current with the latest technology and propose relevant IT solutions.\",\n  },\n  {\n    find: \": 25 February 2025\",\n    replace: \": 03 March 2025\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the 5 textual changes described by the diff:\n#   1. Title: \"Github Foundations Certification Training\"\n#      -> \"Hands-on AI-Assisted Programming Made Simple with GitHub Copilot\"\n#   2. Funding paragraph: course name + framework code update\n#   3. \"Performance gaps\" paragraph rewrite\n#   4. \"Why this course\" paragraph rewrite\n#   5. Date line: \"25 February 2025\" -> \"03 March 2025\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0            # wdFindStop\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$findText, [ref]$true, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$replaceText, [ref]2) | Out-Null\n}\n\n# NOTE: the longer sentence containing \"Github Foundations Certification\n# Training\" as a substring must be replaced BEFORE the standalone title\n# occurrence, otherwise a find/replace-all of the short title also matches\n# (and clobbers) the occurrence embedded inside this sentence.\nReplace-Text `\n    \"We are applying for WSQ funding support for this new course Github Foundations Certification Training according to Software Configuration ICT-DIT-3014-1.1 under Infocomm Technology Framework.\" `\n    \"We are applying for WSQ funding support for this new course Hands-on AI-Assisted Programming Made Simple with GitHub Copilot according to Digital Technology Adoption and Innovation ACC-ICT-3004-1.1 under Infocomm Technology Framework.\"\n\nReplace-Text `\n    \"Github Foundations Certification Training\" `\n    \"Hands-on AI-Assisted Programming Made Simple with GitHub Copilot\"\n\nReplace-Text `\n    \"Many software development teams face challenges in selecting and utilizing the right tools for integration and deployment. A lack of expertise in scripting and automation leads to manual processes that are time-consuming and prone to errors. Moreover, teams often struggle to properly diagnose configuration issues, resulting in prolonged debugging cycles and delayed project timelines.\" `\n    \"One significant challenge is the slow adoption of new technologies and methodologies, hindering the ability to remain competitive. Teams may lack the expertise to effectively integrate emerging tools, which can lead to missed opportunities for automation and improved productivity. Legacy systems and a reluctance to change further compound these issues.\"\n\nReplace-Text `\n    \"The course covers a broad range of Git scripts and tools necessary for effectively integrating and deploying software. Participants will learn how to interpret configuration test results and identify the root causes of issues. This skillset enables them to proactively address problems and implement necessary modifications, which in turn leads to more reliable and efficient software releases.\" `\n    \"This course directly addresses this by providing hands-on experience with cutting-edge AI programming tools. Participants will explore how these tools can streamline organizational coding processes, including using code completion and suggesting code snippets. The course is designed to help individuals stay current with the latest technology and propose relevant IT solutions.\"\n\nReplace-Text \": 25 February 2025\" \": 03 March 2025\"\n"}
